# API testing for parking slot booking
# Updates the sample consumer/provider credential rows and appends a
# parking-booking transcript (consumer side on Sheet1, provider side on
# Sheet2), matching a fresh round of manual API testing.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1 - consumer credentials/booking
$ws2 = $wb.Worksheets.Item(2)   # Sheet2 - provider credentials/booking

# ---------------------------------------------------------------------
# Sheet1 (consumer): swap consumer23 -> consumer41 creds/token
# ---------------------------------------------------------------------
$ws1.Range("A2").Value = "consumer41"
$ws1.Range("C2").Value = "consumer41@gmail.com"
$ws1.Range("F2").Value = "eyJ0eXAiOiJKV1QiLCJhbGciOiJIUzI1NiJ9.eyJ0b2tlbl90eXBlIjoiYWNjZXNzIiwiZXhwIjoxNjUxNTg4MzczLCJqdGkiOiIzMThiMjBhNDIwYjY0OTNlOTViMjRkOGI2MDg2N2UyZiIsInVzZXJfaWQiOjc2OSwicm9sZSI6IkMiLCJ1c2VybmFtZSI6ImNvbnN1bWVyNDEiLCJlbWFpbCI6ImNvbnN1bWVyNDFAZ21haWwuY29tIn0.y7B4-n6c5zo9ETLPsUAyz6b_gh2IN25vJavqKgY-Y3w"

# New booking-confirmation block written under the credentials (rows 4-12)
$ws1.Range("A4").Value = "Mahesh"
$ws1.Range("B4").Value = "Dadeech"
$ws1.Range("C4").NumberFormat = "General"
$ws1.Range("C4").Value = "'8800665544"

$ws1.Range("A5").Value = "TS08UF4343"
$ws1.Range("B5").Value = "'2"
$ws1.Range("C5").Value = "Black SUV"
$ws1.Range("D5").Value = "Mahindra"

$ws1.Range("A6").Value = "'46"

$ws1.Range("A7").Value = "'2022-04-28 12:00"
$ws1.Range("A7").NumberFormat = "m/d/yy h:mm"
$ws1.Range("A7").Font.Color = 0
$ws1.Range("B7").Value = "'2022-05-28 13:00"
$ws1.Range("C7").Value = "'1"
$ws1.Range("D7").Value = "'1"

$ws1.Range("A8").Value = "'43"

$ws1.Range("A10").Value = "Parking ID"
$ws1.Range("B10").Value = "Rating"
$ws1.Range("C10").Value = "Review"

$ws1.Range("A11").Value = "'44"

$ws1.Range("B12").Value = "'3"
$ws1.Range("C12").Value = "nice"

# ---------------------------------------------------------------------
# Sheet2 (provider): swap provider28 -> provider29 creds/token
# ---------------------------------------------------------------------
$ws2.Range("A2").Value = "provider29"
$ws2.Range("C2").Value = "provider29@gmail.com"
$ws2.Range("F2").Value = "eyJ0eXAiOiJKV1QiLCJhbGciOiJIUzI1NiJ9.eyJ0b2tlbl90eXBlIjoiYWNjZXNzIiwiZXhwIjoxNjUxNjM1NjQ3LCJqdGkiOiI0NGQ4NjA1Njk2MWQ0OGM3YWYzNzVlNDRlMmI1ZmYyOCIsInVzZXJfaWQiOjc2OCwicm9sZSI6IlAiLCJ1c2VybmFtZSI6InByb3ZpZGVyMjkiLCJlbWFpbCI6InByb3ZpZGVyMjlAZ21haWwuY29tIn0.8NwxtDSKRdnE0aEWKPX4P73x4oQeXQ4Rafrr8WNEoyI"

# Updated provider phone / pincode used for this booking
$ws2.Range("C4").Value = "'9876342240"
$ws2.Range("D4").Value = "'572143"

# Updated vehicle-entry/exit log row
$ws2.Range("B9").Value = "'44"
$ws2.Range("C9").Value = "TU04KK43"

# ---------------------------------------------------------------------
# View state: Sheet1 becomes the active/selected tab (was Sheet2)
# ---------------------------------------------------------------------
$ws2.Range("C9").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("C14").Select() | Out-Null
